$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order / rename the header row:
#  - C1 was "Order No"  -> now "Order"
#  - D1 was "Active"    -> stays "Active"
#  - E1 was "Category Name" -> now "Category"
# Set E1 first so the new "Category" shared string is created before "Order",
# matching the authoring order used when the columns were rearranged.
$ws.Range("E1").Value = "Category"
$ws.Range("C1").Value = "Order"
$ws.Range("D1").Value = "Active"

# Update the active cell / selection shown when the sheet was last saved.
$ws.Range("E10").Select()
